# Insert a new data row at row 193 (shifting existing rows 193-241 down to
# 194-242) and populate it with a new price observation for
# "Pepino dulce" / Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 193.. down by one, inserting a fresh (blank) row 193.
$ws.Rows.Item(193).Insert()

# Populate the new row 193 with the new observation.
$ws.Cells.Item(193, 1).Value  = 10
$ws.Cells.Item(193, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value  = "La Araucanía"
$ws.Cells.Item(193, 4).Value  = 44736
$ws.Cells.Item(193, 5).Value  = 9
$ws.Cells.Item(193, 6).Value  = 100112043
$ws.Cells.Item(193, 7).Value  = "Pepino dulce"
$ws.Cells.Item(193, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(193, 9).Value  = "Primera"
$ws.Cells.Item(193, 10).Value = 80
$ws.Cells.Item(193, 11).Value = 18000
$ws.Cells.Item(193, 12).Value = 18000
$ws.Cells.Item(193, 13).Value = 18000
$ws.Cells.Item(193, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(193, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(193, 16).Value = 1000
$ws.Cells.Item(193, 17).Value = 18
$ws.Cells.Item(193, 18).Value = "Hortaliza"

# Keep the date column's existing date-number formatting for the new row.
$ws.Cells.Item(193, 4).NumberFormat = $ws.Cells.Item(194, 4).NumberFormat
